# Update the "Login" sheet of the test-data workbook:
#  - refresh the partner-email / partner-password / company-name sample
#    values for the second data row
#  - add a new "PhoneNmbr" column (K) with its sample value
#  - move the active selection over to the newly added column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- updated sample data in the existing columns ---------------------------
$ws.Cells.Item(2, 7).Value  = "digitalmeshlj7ux5@maildrop.cc"   # G2 Partner_Email
$ws.Cells.Item(2, 8).Value  = "DmHsC3Ps!2"                      # H2 Patrner_Password_Record
$ws.Cells.Item(2, 10).Value = "KyPsR"                            # J2 Company_Name

# --- new "PhoneNmbr" column --------------------------------------------------
$ws.Cells.Item(1, 11).Value = "PhoneNmbr"                        # K1 header

# force the phone number to be stored as text (it is a numeric-looking
# string, not a quantity) just like the rest of the sample values
$ws.Cells.Item(2, 11).NumberFormat = "@"
$ws.Cells.Item(2, 11).Value = "3470414515"                       # K2 value

# widen the Partner_Email column slightly (its longest value grew) and give
# the new PhoneNmbr column a sensible best-fit width
$ws.Columns.Item(7).ColumnWidth = 29.166666666666668
$ws.Columns.Item(11).ColumnWidth = 11.333333333333334

# move the selection to the new last column, mirroring the author's click
$ws.Range("K1").Select() | Out-Null
